$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44449, 6, 23, 95.34864439101236),
    @(44450, 0, 20, 82.91186468783683),
    @(44451, 5, 24, 99.4942376254042),
    @(44452, 4, 22, 91.20305115662052),
    @(44453, 5, 22, 91.20305115662052),
    @(44454, 0, 21, 87.05745792222866),
    @(44455, 0, 20, 82.91186468783683),
    @(44456, 4, 18, 74.62067821905315),
    @(44457, 0, 18, 74.62067821905315),
    @(44458, 3, 16, 66.32949175026947),
    @(44459, 4, 16, 66.32949175026947)
)

# Copy the formatting (style) of the last existing date cell (A374) so new date
# cells keep the same look (centered, bold, bordered, date number format).
$ws.Range("A374").Copy() | Out-Null

$r = 375
foreach ($row in $data) {
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $r = $r + 1
}
